# Updates the cryptos list (Price and Volume(1h) columns, plus a couple of
# reordered / replaced coin rows) to match the latest scrape.
#
# Column layout: A = rank (unchanged), B = Coin, C = Link, D = Price, E = Volume(1h)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new text value.
$updates = @{
    2 = @{ D="64.474.72"; E="  -2.07%  " }
    3 = @{ D="3.487.76"; E="  -3.29%  " }
    4 = @{ E="  +0.16%  " }
    5 = @{ D="583.93"; E="  -3.53%  " }
    6 = @{ D="131.10"; E="  -4.62%  " }
    7 = @{ D="3.488.14"; E="  -3.27%  " }
    8 = @{ E="  +0.08%  " }
    9 = @{ D="0.490"; E="  -1.66%  " }
    10 = @{ E="  -1.67%  " }
    11 = @{ D="7.25"; E="  -0.20%  " }
    12 = @{ D="0.386"; E="  -1.96%  " }
    13 = @{ D="4.103.20"; E="  -2.80%  " }
    14 = @{ D="27.59"; E="  -1.78%  " }
    15 = @{ D="0.0000178"; E="  -4.94%  " }
    16 = @{ E="  +0.32%  " }
    17 = @{ D="3.504.40"; E="  -2.79%  " }
    18 = @{ D="64.533.03"; E="  -2.13%  " }
    19 = @{ D="9.87"; E="  -2.57%  " }
    20 = @{ D="14.34"; E="  -2.53%  " }
    21 = @{ D="5.64"; E="  -4.96%  " }
    22 = @{ D="392.45"; E="  -1.52%  " }
    23 = @{ D="0.575"; E="  -2.68%  " }
    24 = @{ D="3.644.50"; E="  -2.86%  " }
    25 = @{ D="73.54"; E="  -1.39%  " }
    26 = @{ E="  -0.01%  " }
    27 = @{ D="0.0000108"; E="  -9.47%  " }
    28 = @{ D="1.56"; E="  -7.56%  " }
    29 = @{ B="Binance-PegBSC-USD"; C="https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D="1.00"; E="  +0.04%  " }
    30 = @{ B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="7.40"; E="  -10.10%  " }
    31 = @{ D="2.25"; E="  -6.52%  " }
    32 = @{ D="8.14"; E="  -6.19%  " }
    33 = @{ D="3.504.24"; E="  -2.86%  " }
    34 = @{ E="  +0.03%  " }
    35 = @{ D="24.00"; E="  -2.71%  " }
    36 = @{ D="0.146"; E="  -1.78%  " }
    37 = @{ D="1.57"; E="  -2.60%  " }
    38 = @{ D="5.20"; E="  -3.63%  " }
    39 = @{ D="170.81"; E="  -0.77%  " }
    40 = @{ D="6.93"; E="  -2.53%  " }
    41 = @{ D="0.0802"; E="  -4.27%  " }
    42 = @{ D="0.812"; E="  -3.58%  " }
    43 = @{ D="26.18"; E="  +0.09%  " }
    44 = @{ E="  +0.08%  " }
    45 = @{ D="41.88"; E="  -3.60%  " }
    46 = @{ D="1.20"; E="  -4.56%  " }
    47 = @{ D="4.34"; E="  -4.68%  " }
    48 = @{ D="1.63"; E="  -4.76%  " }
    49 = @{ B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="2.433.49"; E="  -1.16%  " }
    50 = @{ B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="6.85"; E="  -3.54%  " }
    51 = @{ B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0266"; E="  -2.65%  " }

}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $value = $rowData[$col]
        $cell = $ws.Range("$col$rowNum")

        if ($col -eq "D") {
            # Price values such as "7.25" or "1.00" would otherwise be
            # auto-recognized as numbers by Excel; prefix with an apostrophe
            # (quote-prefix) so they are stored as literal text, matching the
            # original sheet where every Price cell is text. Values that
            # already contain more than one "." (e.g. "64.474.72") are never
            # parsed as numbers, so they don't need the prefix.
            $dotCount = ($value.ToCharArray() | Where-Object { $_ -eq '.' } | Measure-Object).Count
            if ($dotCount -eq 1) {
                $cell.Value = "'" + $value
            } else {
                $cell.Value = $value
            }
        } else {
            $cell.Value = $value
        }
    }
}
